$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.596.59'
$ws.Range("E2").Value = '  -2.42%  '
$ws.Range("D3").Value = '2.893.90'
$ws.Range("E3").Value = '  -2.07%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '''569.08'
$ws.Range("E5").Value = '  -4.34%  '
$ws.Range("D6").Value = '''144.20'
$ws.Range("E6").Value = '  -3.02%  '
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("E8").Value = '  -0.39%  '
$ws.Range("D9").Value = '2.895.07'
$ws.Range("E9").Value = '  -1.96%  '
$ws.Range("D10").Value = '''6.99'
$ws.Range("E10").Value = '  -3.91%  '
$ws.Range("D11").Value = '''0.146'
$ws.Range("E11").Value = '  -2.92%  '
$ws.Range("E12").Value = '  -2.23%  '
$ws.Range("E13").Value = '  -0.99%  '
$ws.Range("D14").Value = '''32.02'
$ws.Range("E14").Value = '  -2.51%  '
$ws.Range("E15").Value = '  -0.63%  '
$ws.Range("D16").Value = '3.373.43'
$ws.Range("E16").Value = '  -2.02%  '
$ws.Range("D17").Value = '61.604.38'
$ws.Range("D18").Value = '''6.56'
$ws.Range("E18").Value = '  -2.06%  '
$ws.Range("D19").Value = '2.892.75'
$ws.Range("E19").Value = '  -2.03%  '
$ws.Range("D20").Value = '''433.32'
$ws.Range("E20").Value = '  -2.29%  '
$ws.Range("D21").Value = '''13.11'
$ws.Range("E21").Value = '  -2.81%  '
$ws.Range("D22").Value = '''0.657'
$ws.Range("E22").Value = '  -1.58%  '
$ws.Range("E23").Value = '  -2.75%  '
$ws.Range("D24").Value = '''79.36'
$ws.Range("E24").Value = '  -2.13%  '
$ws.Range("D25").Value = '''12.03'
$ws.Range("E25").Value = '  +2.24%  '
$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").Value = '''1.00'
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("B27").Value = 'RenderToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D27").Value = '''10.01'
$ws.Range("E27").Value = '  -11.06%  '
$ws.Range("E28").Value = '  -5.50%  '
$ws.Range("E29").Value = '  +2.71%  '
$ws.Range("E30").Value = '  -3.40%  '
$ws.Range("D31").Value = '''2.50'
$ws.Range("E31").Value = '  -4.17%  '
$ws.Range("D32").Value = '''2.05'
$ws.Range("E32").Value = '  -7.28%  '
$ws.Range("E33").Value = '  +0.02%  '
$ws.Range("D34").Value = '''0.107'
$ws.Range("E34").Value = '  -1.89%  '
$ws.Range("D35").Value = '''25.51'
$ws.Range("E35").Value = '  -3.64%  '
$ws.Range("D36").Value = '''0.960'
$ws.Range("E36").Value = '  -3.20%  '
$ws.Range("D37").Value = '''5.40'
$ws.Range("E37").Value = '  -3.38%  '
$ws.Range("D38").Value = '''48.90'
$ws.Range("E38").Value = '  -1.62%  '
$ws.Range("E39").Value = '  -5.78%  '
$ws.Range("E40").Value = '  -10.13%  '
$ws.Range("E41").Value = '  -2.78%  '
$ws.Range("E42").Value = '  -3.32%  '
$ws.Range("D43").Value = '''39.58'
$ws.Range("E43").Value = '  +1.74%  '
$ws.Range("E44").Value = '  -5.43%  '
$ws.Range("D45").Value = '2.706.37'
$ws.Range("E45").Value = '  +0.46%  '
$ws.Range("D46").Value = '''132.47'
$ws.Range("E46").Value = '  -2.30%  '
$ws.Range("E47").Value = '  -0.51%  '
$ws.Range("D48").Value = '''347.51'
$ws.Range("E48").Value = '  -3.63%  '
$ws.Range("E50").Value = '  -1.22%  '
$ws.Range("D51").Value = '''21.64'
$ws.Range("E51").Value = '  -5.39%  '
